# Apply the change: swap the contents of columns C and D (OR_LowerCI / OR_UpperCI)
# for the rows whose values are computed via formulas on sheet "PO_reduced",
# then move the active cell selection to D19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PO_reduced")

# Rows that hold formulas in columns C and D and need to be swapped.
$rows = @(3, 8, 9, 10, 11, 14, 16)

foreach ($r in $rows) {
    $cCell = $ws.Cells.Item($r, 3)   # column C
    $dCell = $ws.Cells.Item($r, 4)   # column D

    $cFormula = $cCell.Formula
    $dFormula = $dCell.Formula

    $cCell.Formula = $dFormula
    $dCell.Formula = $cFormula
}

# Activate the sheet and update the selected cell, matching the diff.
$ws.Activate()
$ws.Range("D19").Select()
